# Update the "Förändrad" (Changed) date column (C) for rows 2-8
# from 2023-09-01 (serial 45170) to 2023-09-05 (serial 45174).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45174
}
